# 8/30 update Tyee data and fisheries openings
# Adds two new fisheries-opening rows to the "Skeena" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skeena")

# Row 38 - FN0877
$ws.Range("A38").Value = "FN0877"
$ws.Range("B38").Value = "Aboriginal"
$ws.Range("C38").Value = "Aug 29-31"
$ws.Range("D38").Value = "Sockeye"
$ws.Range("E38").Value = "Selective Gear"
$ws.Range("F38").Value = "Region 6-Lake Babine Nation"
$ws.Range("G38").Value = 3
$ws.Range("I38").Value = "Sockeye target, selective gear only"

# Row 39 - FN0883
$ws.Range("A39").Value = "FN0883"
$ws.Range("B39").Value = "Aboriginal"
$ws.Range("C39").Value = "Aug 30-Sept 5"
$ws.Range("D39").Value = "Sockeye "
$ws.Range("E39").Value = "Selective Gear"
$ws.Range("F39").Value = "Region 6-Lake Babine Nation"
$ws.Range("G39").Value = 7
$ws.Range("I39").Value = "Sockeye target, selective gear only"

# Match the date-column format used by the rest of the table (column C)
$ws.Range("C38").NumberFormat = $ws.Range("C37").NumberFormat
$ws.Range("C39").NumberFormat = $ws.Range("C37").NumberFormat
